$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Continue the participant list: rows 2-19 already hold "name 1".."name 18".
# Append "name 19" through "name 30" in rows 20-31.
$names = @(
    "name 19", "name 20", "name 21", "name 22", "name 23", "name 24",
    "name 25", "name 26", "name 27", "name 28", "name 29", "name 30"
)

$row = 20
foreach ($name in $names) {
    $ws.Cells.Item($row, 1).Value = $name
    $row = $row + 1
}

# Match the author's final selection: A19 active cell, A19:A31 selected.
$ws.Range("A19:A31").Select()
